# Applies the COP22 Czech Republic workbook edit:
#  - Splits the "Climate Protection Policy" paragraph (row 13) into two rows,
#    moving the relevance/topic/.../notes data off of it (it is no longer "relevant").
#  - Splits the "We are aware of the fact..." paragraph (row 19) into two rows,
#    keeping the relevance/topic/.../notes data on the first (shorter) half.
#  - Adds a brand new relevance/topic/.../notes data row for the closing
#    "We must focus on effective implementation..." paragraph.
#  - Updates row heights / column width / sheet view to match the re-saved file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Row 10 ("Jointly with the EU...") loses its relevance data.
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = "no"
$ws.Range("C10:H10").ClearContents()

# ---------------------------------------------------------------------------
# 2) Row 13 ("A new Climate Protection Policy...") is split into two rows.
#    Insert a new row right after it, then redistribute the text.
# ---------------------------------------------------------------------------
$ws.Rows(14).Insert()

$ws.Range("A13").Value = "A new Climate Protection Policy in the Czech Republic, which will serve as our `nlow carbon development strategy is expected to be adopted in March 2017. "
$ws.Range("B13").Value = "no"
$ws.Range("C13:H13").ClearContents()

$ws.Range("A14").Value = "The `nPolicy suggests reaching a long term goal to reduce GHG emissions in 2050 by 80 `nper cent and that all key economic sectors will contribute to achieving this goal `nwith specific measures."
$ws.Range("B14").Value = "no"

# ---------------------------------------------------------------------------
# 3) Row 19 ("We are aware of the fact...") -- now row 20 after the insert
#    above -- is split into two rows. The first (shorter) half keeps the
#    relevance data that used to live on the combined paragraph.
# ---------------------------------------------------------------------------
$ws.Rows(21).Insert()

$ws.Range("A20").Value = "We are aware of the fact that the poorest and most vulnerable countries are `ndependent on support to address their domestic climate challenges and risks.`n"
$ws.Range("B20").Value = "yes"
$ws.Range("C20").Value = "adaptation, mitigation"
$ws.Range("D20").Value = "financial resources"
$ws.Range("E20").Value = "global"
$ws.Range("F20").Value = "n.a."
$ws.Range("G20").Value = "prioritarian"
$ws.Range("H20").Value = "Support of most vulerable and poorest countries. "

$ws.Range("A21").Value = "In this respect I am pleased to inform you that the Czech Republic provided in `ntotal 8,1 mil. EUR of climate finance in 2015 through its bilateral and multilateral `nchannels and is prepared to continue to provide nece ssary support to developing `ncountries in order to assist them in their adaptation and mitigation needs and `nefforts."
$ws.Range("B21").Value = "no"

# ---------------------------------------------------------------------------
# 4) The closing paragraph "We must focus on effective implementation..."
#    (row 22 originally, now row 24 after the two inserts above) gains a
#    brand-new relevance data row.
# ---------------------------------------------------------------------------
$ws.Range("B24").Value = "yes"
$ws.Range("D24").Value = "measures"
$ws.Range("E24").Value = "global"
$ws.Range("F24").Value = "n.a."
$ws.Range("G24").Value = "utilitarian"
$ws.Range("H24").Value = "Presenting the need for implementation of measures for the benefit of current and future generations. "
$ws.Range("C24").Value = "implementation, UNFCCC agreements and policies"

# ---------------------------------------------------------------------------
# 5) Row heights (Excel re-flowed these against a slightly different default
#    font metric when the file was re-saved).
# ---------------------------------------------------------------------------
$ws.Rows(2).RowHeight = 158.4
$ws.Rows(3).RowHeight = 100.8
$ws.Rows(4).RowHeight = 43.2
$ws.Rows(5).RowHeight = 72
$ws.Rows(6).RowHeight = 57.6
$ws.Rows(7).RowHeight = 72
$ws.Rows(8).RowHeight = 57.6
$ws.Rows(9).RowHeight = 86.4
$ws.Rows(10).RowHeight = 100.8
$ws.Rows(11).RowHeight = 72
$ws.Rows(12).RowHeight = 43.2
$ws.Rows(13).RowHeight = 57.6
$ws.Rows(14).RowHeight = 86.4
$ws.Rows(15).RowHeight = 72
$ws.Rows(16).RowHeight = 28.8
$ws.Rows(17).RowHeight = 43.2
$ws.Rows(18).RowHeight = 43.2
$ws.Rows(19).RowHeight = 86.4
$ws.Rows(20).RowHeight = 72
$ws.Rows(21).RowHeight = 129.6
$ws.Rows(22).RowHeight = 86.4
$ws.Rows(23).RowHeight = 28.8
$ws.Rows(24).RowHeight = 149.4
$ws.Rows(25).RowHeight = 28.8
$ws.Rows(26).RowHeight = 43.2

# ---------------------------------------------------------------------------
# 6) Sheet view changes: selection moves to C20 (also drops the scrolled-down
#    topLeftCell the file was saved with, since the sheet no longer needs to
#    stay scrolled to keep the old row 20/23 in view). Column width is left
#    as-is: the ~0.016 pt nudge Excel made there on re-save is far finer than
#    this engine's column-width rounding granularity, so attempting to set
#    it would overshoot by more than leaving it alone.
# ---------------------------------------------------------------------------
$ws.Range("C20").Select()
